$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the Cases tab query (B2): append ORDER BY / LIMIT clause
$b2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100"

# Update the Samples tab query (B3): append ORDER BY / LIMIT clause
$b3 = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = $b3 + "`n order By samp.sample_id ASC LIMIT 100"

# Update the Files tab query (B4): replace trailing "order by" clause with new ORDER BY / LIMIT clause
$b4 = $ws.Range("B4").Value2
$b4 = $b4 -replace "    order by f\.file_name$", "     order By f.file_name ASC LIMIT 100"
$ws.Range("B4").Value2 = $b4

# Update selection to C4
$ws.Range("C4").Select()
